# Insert a brand-new weekly observation row before current row 989, which
# shifts the existing rows 989..1028 down to 990..1029 (dimension grows to
# A1:R1029). Then populate the new row 989 with the new data point, keeping
# the constant/categorical columns identical to the (now shifted) row below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 989, pushing everything at/after it down by one.
$ws.Rows.Item(989).EntireRow.Insert()

# Populate the newly inserted row 989 with the new data point.
$ws.Cells.Item(989, 1).Value2 = 6
$ws.Cells.Item(989, 2).Value2 = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(989, 3).Value2 = 'Metropolitana'
$ws.Cells.Item(989, 4).Value2 = 44826
$ws.Cells.Item(989, 5).Value2 = 13
$ws.Cells.Item(989, 6).Value2 = 100112028
$ws.Cells.Item(989, 7).Value2 = 'Sandia'
$ws.Cells.Item(989, 8).Value2 = 'Sin especificar'
$ws.Cells.Item(989, 9).Value2 = 'Primera'
$ws.Cells.Item(989, 10).Value2 = 4400
$ws.Cells.Item(989, 11).Value2 = 850
$ws.Cells.Item(989, 12).Value2 = 900
$ws.Cells.Item(989, 13).Value2 = 872
$ws.Cells.Item(989, 14).Value2 = '$/kilo (volumen en unidades)'
$ws.Cells.Item(989, 15).Value2 = 'Perú'
$ws.Cells.Item(989, 16).Value2 = 872
$ws.Cells.Item(989, 17).Value2 = 1
$ws.Cells.Item(989, 18).Value2 = 'Hortaliza'
